# Add the two newly sold special-price items described in the commit
# message: "La passion de Dodin Buffant" (06.10.24) and "Garfield" (08.09.24)
# The Garfield (08.09.24 / serial 45543) entries already exist in the sheet
# (rows 16-17); the still-missing entry is the new sales row for
# "La passion de Dodin Buffant" on 06.10.2024 (serial date 45571).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended right after the current last row (18)
$newRow = 19

$ws.Cells.Item($newRow, 1).Value = 45571          # Datum: 06.10.2024
$ws.Cells.Item($newRow, 2).Value = "Spez 1"       # Spezialpreis
$ws.Cells.Item($newRow, 3).Value = "Kaffee und Gipfeli"  # Artikelname
$ws.Cells.Item($newRow, 4).Value = 5              # Verkaufspreis
$ws.Cells.Item($newRow, 5).Value = 6              # Anzahl verkaufter Artikel

# Extend Table1 so the new row becomes part of the table range (A1:E18 -> A1:E19)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E" + $newRow))

# Match the saved selection state from the edit (active cell H19)
$ws.Range("H19").Select() | Out-Null
